# The documentation table tracks, per feature, whether it has been
# implemented in the API / REST API / WEBAPP layers.
#
# 1) "Dodaj ucznia" (Add student) row, REST API column: the cell already
#    contained "X " but split across two separate runs (<w:t>X</w:t> and
#    <w:t xml:space="preserve"> </w:t>). Re-running Find/Replace over the
#    exact same text collapses it back down into a single merged run,
#    matching how Word normally stores text after an in-place edit.
#
# 2) "Pobierz uczniow po stringu (...)" (Get students by string) row:
#    mark the feature as implemented by typing "X" into the, until now
#    empty, API and REST API columns.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1) merge the "X " runs in the "Dodaj ucznia" row / REST API column ---
$cellAddStudentRest = $t.Cell(2, 3)
[void]$cellAddStudentRest.Range.Find.Execute("X ", $false, $false, $false, `
    $false, $false, $true, 1, $false, "X ", 2)

# --- 2) mark "Pobierz uczniow po stringu" as done for API and REST API ---
$rowQueryByString = 12

$apiCell = $t.Cell($rowQueryByString, 2)
$apiRange = $apiCell.Range
$apiRange.InsertAfter("X")
$apiCell = $t.Cell($rowQueryByString, 2)
$apiCell.Range.Font.Size = 12
$apiCell.Range.Font.SizeBi = 12

$restCell = $t.Cell($rowQueryByString, 3)
$restRange = $restCell.Range
$restRange.InsertAfter("X")
$restCell = $t.Cell($rowQueryByString, 3)
$restCell.Range.Font.Size = 12
$restCell.Range.Font.SizeBi = 12
